$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") changes from 45183 to 45184 for every data row (2..62)
for ($r = 2; $r -le 62; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value = 45184
    }
}

# Row 2 hyperlink formulas gain a friendly-name second argument.
# S2 reproduces the source data's malformed quoting verbatim (missing
# closing quote before the comma, stray trailing quote).
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/artfynd/A 31987-2023.xlsx, "A 31987-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/kartor/A 31987-2023.png", "A 31987-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/klagomål/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/klagomålsmail/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/tillsyn/A 31987-2023.docx", "A 31987-2023")'

# Y2 was stored as a literal inline string (not a real formula); it becomes
# a genuine HYPERLINK formula like its siblings.
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_ESKILSTUNA/tillsynsmail/A 31987-2023.docx", "A 31987-2023")'
